$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (G and H) for "start date" and "end date" of an issue.
# This shifts the existing Summary / Parent issue / Parent summary columns
# from G,H,I to I,J,K.
$ws.Columns("G:H").Insert()

# Fill in the new start/end date columns for the three data rows that have
# time entries (row 4 is a "parent issue only" row with no dates).
$ws.Range("G1:G3").Value = 43840
$ws.Range("H1:H3").Value = 43845

$ws.Range("G1:H3").NumberFormat = "yyyy-mm-dd"

# Match the column width used for the new date columns (~11.17 Excel
# "characters", expressed here in the pre-padding input unit this engine
# expects for ColumnWidth).
$ws.Columns("G:H").ColumnWidth = 10.336666666666666

# Update the active selection to reflect where the user ended up.
$null = $ws.Range("H2").Select()
